$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value into a cell while preserving its original
# (default/General) style, even when the text looks like a number
# ("0.9947", "1.000", etc). We briefly force Text format so Excel
# does not coerce the string to a numeric value, then clear the
# explicit formatting we just applied so the cell style reverts to
# its original (unstyled / style index 0) state.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '26.351.17'
$ws.Range("E2").Value = '  -4.18%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.768.08'
$ws.Range("E3").Value = '  -3.07%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.9947'
$ws.Range("E4").Value = '  -0.65%  '

# Row 5
Set-TextValue $ws.Range("D5") '317.42'
$ws.Range("E5").Value = '  +1.58%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.9954'
$ws.Range("E6").Value = '  -0.56%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.4274'
$ws.Range("E7").Value = '  +0.85%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3638'
$ws.Range("E8").Value = '  +1.16%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.07171'
$ws.Range("E9").Value = '  -0.23%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.8484'
$ws.Range("E10").Value = '  -0.68%  '

# Row 11
Set-TextValue $ws.Range("D11") '20.34'
$ws.Range("E11").Value = '  -0.34%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.772.83'
$ws.Range("E12").Value = '  -5.35%  '

# Row 13
Set-TextValue $ws.Range("D13") '5.241'
$ws.Range("E13").Value = '  -2.14%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.424'
$ws.Range("E14").Value = '  +0.08%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.06842'
$ws.Range("E15").Value = '  -1.00%  '

# Row 16
Set-TextValue $ws.Range("D16") '1.000'
$ws.Range("E16").Value = '  -0.27%  '

# Row 17
Set-TextValue $ws.Range("D17") '78.75'
$ws.Range("E17").Value = '  -3.74%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.000008669'
$ws.Range("E18").Value = '  -2.45%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.9968'
$ws.Range("E19").Value = '  -0.73%  '

# Row 20
Set-TextValue $ws.Range("D20") '14.96'
$ws.Range("E20").Value = '  -2.32%  '

# Row 21
Set-TextValue $ws.Range("D21") '26.355.93'
$ws.Range("E21").Value = '  -4.77%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.091'
$ws.Range("E22").Value = '  -0.18%  '

# Row 23
$ws.Range("E23").Value = '  +0.76%  '

# Row 24
Set-TextValue $ws.Range("D24") '1.982.48'
$ws.Range("E24").Value = '  -4.29%  '

# Row 25
Set-TextValue $ws.Range("D25") '151.52'
$ws.Range("E25").Value = '  -1.56%  '

# Row 26
$ws.Range("E26").Value = '  -6.15%  '

# Row 27
$ws.Range("E27").Value = '  -1.64%  '

# Row 28
Set-TextValue $ws.Range("D28") '5.093'
$ws.Range("E28").Value = '  -0.85%  '

# Row 29
Set-TextValue $ws.Range("D29") '113.66'
$ws.Range("E29").Value = '  -0.17%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.802'
$ws.Range("E30").Value = '  +2.50%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.08933'
$ws.Range("E31").Value = '  +0.28%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.7284'
$ws.Range("E32").Value = '  -2.45%  '

# Row 33
Set-TextValue $ws.Range("D33") '1.124'
$ws.Range("E33").Value = '  +0.28%  '

# Row 34
Set-TextValue $ws.Range("D34") '4.315'
$ws.Range("E34").Value = '  -4.32%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.9949'
$ws.Range("E35").Value = '  -0.65%  '

# Row 36
Set-TextValue $ws.Range("D36") '2.718'
$ws.Range("E36").Value = '  -7.92%  '

# Row 37
Set-TextValue $ws.Range("D37") '1.090'
$ws.Range("E37").Value = '  +0.93%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.05151'
$ws.Range("E38").Value = '  -1.47%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.01891'
$ws.Range("E39").Value = '  -1.11%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.4922'
$ws.Range("E40").Value = '  -2.35%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.1609'
$ws.Range("E41").Value = '  -2.45%  '

# Row 42
Set-TextValue $ws.Range("D42") '2.579'
$ws.Range("E42").Value = '  -6.95%  '

# Row 43
Set-TextValue $ws.Range("D43") '6.278'
$ws.Range("E43").Value = '  -0.68%  '

# Row 44
Set-TextValue $ws.Range("D44") '8.014'
$ws.Range("E44").Value = '  -3.64%  '

# Row 45
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D45") '104.78'
$ws.Range("E45").Value = '  -1.18%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D46") '10.23'
$ws.Range("E46").Value = '  -1.51%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.9946'
$ws.Range("E47").Value = '  -0.60%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.4485'
$ws.Range("E48").Value = '  -3.43%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.06180'
$ws.Range("E49").Value = '  -4.04%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.610'
$ws.Range("E50").Value = '  -0.06%  '

# Row 51
Set-TextValue $ws.Range("D51") '1.738'
$ws.Range("E51").Value = '  +2.79%  '
